$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new literal text value, per the commit diff (Price/Volume columns).
$changes = @{
    "D2" = "261.26"
    "E2" = "0.89%"
    "E3" = "0.70%"
    "D4" = "4.709"
    "E4" = "0.68%"
    "D5" = "0.06209"
    "E5" = "2.84%"
    "D6" = "6.729"
    "E6" = "0.81%"
    "D7" = "0.8506"
    "E7" = "-0.96%"
    "D8" = "0.9076"
    "E8" = "-1.01%"
    "D9" = "0.1405"
    "E9" = "0.65%"
    "D10" = "0.04695"
    "E10" = "-12.75%"
    "D11" = "0.07099"
    "E11" = "0.25%"
    "D12" = "0.03173"
    "E12" = "3.18%"
    "D13" = "0.09056"
    "E13" = "-0.81%"
    "D14" = "0.001533"
    "E14" = "0.03%"
    "D15" = "0.0006182"
    "E15" = "2.32%"
    "D16" = "0.006152"
    "E16" = "0.25%"
    "E17" = "0.04%"
    "E18" = "-0.15%"
    "E19" = "0.58%"
    "E20" = "-0.71%"
    "E21" = "0.21%"
    "D22" = "4.107"
    "E22" = "-1.09%"
    "D23" = "0.04227"
    "E23" = "-0.09%"
    "E24" = "-0.04%"
    "D25" = "0.004130"
    "E25" = "2.59%"
    "E26" = "0.12%"
    "E27" = "6.13%"
    "D40" = "0.03918"
    "E40" = "1.82%"
    "D41" = "0.1113"
    "E41" = "-0.12%"
    "D42" = "0.004133"
    "E42" = "2.71%"
    "E43" = "-0.70%"
    "D44" = "0.01390"
    "E44" = "-8.14%"
    "D45" = "0.00005174"
    "E45" = "1.11%"
    "E46" = "0.12%"
    "D47" = "0.03591"
    "E47" = "-34.15%"
    "E48" = "26.80%"
    "D49" = "0.00002102"
    "E49" = "0.12%"
    "D50" = "0.0002002"
    "E50" = "0.12%"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    # Leading apostrophe forces Excel to store the value as literal text,
    # matching the workbook's existing inline-string (text) cells exactly
    # -- these are tickers/prices/percentages stored as text, not numbers.
    $cell.Value = "'" + $changes[$addr]
    $cell.Style = $origStyle
}
